$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.321.27'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.136.37'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.11%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '609.64'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.88'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.11%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.133.13'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.15%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.530'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.15%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.44%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.34'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.65%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.474'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.52%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000256'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.46%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.40'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.18%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.647.49'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.20%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.117'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.77%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.229.53'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.00%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.133.22'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.84'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '476.77'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.63'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.720'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.79'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.79%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.53'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.70%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.98'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.88%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.78'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.65%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.50'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.77%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.35'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +9.91%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +3.90%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.08'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -4.49%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.65'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.42%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.62'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.52%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.10'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.21%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.91'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.67%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '52.49'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.96%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0731'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +3.42%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '450.80'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.68%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.97'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +3.51%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.37%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.11%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.31'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.96%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.850.90'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.266'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.32%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.25'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.16%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +5.70%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '26.26'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.114'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.34%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '120.06'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.37%  '
